# Reorder/update the NBA roster table on the active worksheet (A2:C17)
# to reflect the new player/position/team data from the published update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Shai Gilgeous-Alexander", "PG", "Oklahoma City Thunder"),
    @("T.J. McConnell", "PG", "Indiana Pacers"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Mike Conley", "PG", "Minnesota Timberwolves"),
    @("Dennis Schröder", "PG", "Brooklyn Nets"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("RJ Barrett", "SF,PF", "Toronto Raptors"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons")
)

$row = 2
foreach ($record in $data) {
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $row++
}
